# Update cryptocurrency price/volume figures per the Wed Aug 23 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.062.78'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.635.38'
$ws.Range("E3").Value = '  -1.91%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''213.08'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''0.2590'
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Value = '''0.06285'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '''20.68'
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("D11").Value = '''0.07657'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").Value = '1.646.13'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("D13").Value = '''4.410'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '1.859.16'
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").Value = '''0.5500'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '0.0₅8163'
$ws.Range("D17").Value = '''64.85'
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("D18").Value = '26.052.66'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '''4.683'
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").Value = '''188.15'
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").Value = '''10.14'
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").Value = '''6.138'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -2.77%  '
$ws.Range("E26").Value = '  -2.84%  '
$ws.Range("D27").Value = '''7.390'
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("D29").Value = '''1.398'
$ws.Range("E29").Value = '  +3.11%  '
$ws.Range("D30").Value = '''0.05948'
$ws.Range("E30").Value = '  -5.04%  '
$ws.Range("D31").Value = '''1.254'
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("D32").Value = '''3.437'
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("D33").Value = '''3.401'
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").Value = '''1.633'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").Value = '''0.9822'
$ws.Range("E35").Value = '  -1.61%  '
$ws.Range("D37").Value = '''2.762'
$ws.Range("E37").Value = '  +1.17%  '
$ws.Range("D38").Value = '''0.5704'
$ws.Range("E38").Value = '  -5.43%  '
$ws.Range("D39").Value = '''0.01614'
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").Value = '''0.8512'
$ws.Range("E40").Value = '  -2.37%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").Value = '''5.734'
$ws.Range("E42").Value = '  -6.21%  '
$ws.Range("D43").Value = '1.032.54'
$ws.Range("E43").Value = '  -6.69%  '
$ws.Range("D44").Value = '''100.29'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").Value = '1.785.63'
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").Value = '''55.82'
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("E47").Value = '  -6.35%  '
$ws.Range("D48").Value = '''1.0000'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = '''8.015'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("E50").Value = '  -1.41%  '
$ws.Range("E51").Value = '  -0.59%  '
